$d = $word.ActiveDocument

# The first-page date is in its own paragraph (centered, large font).
# Replace only that standalone occurrence of "2017-04-10" with "2017-04-24",
# leaving the other occurrence ("2017-04-10 - 2017-06-01" in the schedule
# table) untouched.
$para = $d.Paragraphs(3)
$range = $para.Range
$range.Find.Execute("2017-04-10", $true, $false, $false, $false, $false,
                     $true, 1, $false, "2017-04-24", 2)
